$wb = $excel.ActiveWorkbook

# ALC row 10
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 704
$ws.Range("I10").Value = 704
$ws.Range("K10").Value = 704
$ws.Range("M10").Value = -411

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 243.59091
$ws.Range("I33").Value = 270.21054
$ws.Range("J33").Value = 75
$ws.Range("K33").Value = 270.21054
$ws.Range("L33").Value = 75
$ws.Range("M33").Value = -41.21053999999998
$ws.Range("N33").Value = -533

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1352.5714
$ws.Range("I40").Value = 1395
$ws.Range("J40").Value = 1296
$ws.Range("K40").Value = 1395
$ws.Range("L40").Value = 1296
$ws.Range("M40").Value = -1220
$ws.Range("N40").Value = -1646

# ALC row 44
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

# ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 55251
$ws.Range("I86").Value = 2000
$ws.Range("J86").Value = 73001.336
$ws.Range("K86").Value = 2000
$ws.Range("L86").Value = 73001.336
$ws.Range("M86").Value = -877
$ws.Range("N86").Value = -75247.336

# ALC row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 55251
$ws.Range("I89").Value = 2000
$ws.Range("J89").Value = 73001.336
$ws.Range("K89").Value = 10000
$ws.Range("L89").Value = 365006.68
$ws.Range("M89").Value = -4384
$ws.Range("N89").Value = -376238.68

# ALC row 94
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 461
$ws.Range("I94").Value = 461
$ws.Range("K94").Value = 461
$ws.Range("M94").Value = -10

# ALC row 96
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 9057.143
$ws.Range("I96").Value = 13550
$ws.Range("J96").Value = 5687.5
$ws.Range("K96").Value = 40650
$ws.Range("L96").Value = 17062.5
$ws.Range("M96").Value = -39277
$ws.Range("N96").Value = -19808.5

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2464.2856
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 2464.2856
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 7392.8568
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -9608.856800000001

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1520.625
$ws.Range("I132").Value = 1520.625
$ws.Range("K132").Value = 4561.875
$ws.Range("M132").Value = -2031.875

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 880.7143
$ws.Range("I2").Value = 794.1667
$ws.Range("K2").Value = 794.1667
$ws.Range("M2").Value = -681.1667

# ARM row 30
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 450
$ws.Range("I30").Value = 400
$ws.Range("J30").Value = 500
$ws.Range("K30").Value = 400
$ws.Range("L30").Value = 500
$ws.Range("M30").Value = -250
$ws.Range("N30").Value = -800

# ARM row 33
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 4026
$ws.Range("I33").Value = 4026
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 4026
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -3697
$ws.Range("N33").ClearContents()

# ARM row 36
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

# ARM row 39
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 6766.125
$ws.Range("I39").Value = 5872
$ws.Range("J39").Value = 7302.6
$ws.Range("K39").Value = 5872
$ws.Range("L39").Value = 7302.6
$ws.Range("M39").Value = -5352
$ws.Range("N39").Value = -8342.6

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1299.4445
$ws.Range("I97").Value = 1142
$ws.Range("J97").Value = 1708.8
$ws.Range("K97").Value = 1142
$ws.Range("L97").Value = 1708.8
$ws.Range("M97").Value = -646
$ws.Range("N97").Value = -2700.8

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 880.7143
$ws.Range("I116").Value = 794.1667
$ws.Range("K116").Value = 794.1667
$ws.Range("M116").Value = 1499.8333

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 9985
$ws.Range("J122").Value = 8462.5
$ws.Range("L122").Value = 25387.5
$ws.Range("N122").Value = -30287.5

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 880.7143
$ws.Range("I3").Value = 794.1667
$ws.Range("K3").Value = 794.1667
$ws.Range("M3").Value = -680.1667

# BSM row 36
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 6349.6
$ws.Range("I36").Value = 6349.6
$ws.Range("K36").Value = 6349.6
$ws.Range("M36").Value = -5815.6

# BSM row 38
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 51250
$ws.Range("J38").Value = 55000
$ws.Range("L38").Value = 55000
$ws.Range("N38").Value = -55832

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 809.8182
$ws.Range("I22").Value = 800.7619
$ws.Range("K22").Value = 800.7619
$ws.Range("M22").Value = -450.7619

# CRP row 33
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 38333.668
$ws.Range("I33").Value = 27500.5
$ws.Range("K33").Value = 27500.5
$ws.Range("M33").Value = -27121.5

# CRP row 38
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 8012.6665
$ws.Range("I38").Value = 8012.6665
$ws.Range("K38").Value = 8012.6665
$ws.Range("M38").Value = -7635.6665

# CRP row 46
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H46").Value = 8012.6665
$ws.Range("I46").Value = 8012.6665
$ws.Range("K46").Value = 8012.6665
$ws.Range("M46").Value = -7801.6665

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 716999.3
$ws.Range("I99").Value = 2873.75
$ws.Range("K99").Value = 2873.75
$ws.Range("M99").Value = -1375.75

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 716999.3
$ws.Range("I126").Value = 2873.75
$ws.Range("K126").Value = 8621.25
$ws.Range("M126").Value = -6151.25

# CUL row 9
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 5588.25
$ws.Range("I9").Value = 3300
$ws.Range("J9").Value = 6351
$ws.Range("K9").Value = 9900
$ws.Range("L9").Value = 19053
$ws.Range("M9").Value = -9676
$ws.Range("N9").Value = -19501

# CUL row 46
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 3423.3333
$ws.Range("I46").Value = 120
$ws.Range("K46").Value = 360
$ws.Range("M46").Value = -269

# CUL row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1522.75
$ws.Range("I129").Value = 1029
$ws.Range("K129").Value = 3087
$ws.Range("M129").Value = 1913

# GSM row 57
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 37749.25
$ws.Range("J57").Value = 37749.25
$ws.Range("L57").Value = 37749.25
$ws.Range("N57").Value = -39389.25

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1180
$ws.Range("J80").Value = 1180
$ws.Range("L80").Value = 1180
$ws.Range("N80").Value = -3176

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 1180
$ws.Range("J83").Value = 1180
$ws.Range("L83").Value = 5900
$ws.Range("N83").Value = -15884

# LTW row 29
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3500
$ws.Range("I62").Value = 2000
$ws.Range("K62").Value = 2000
$ws.Range("M62").Value = -1376

# WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 3500
$ws.Range("I65").Value = 2000
$ws.Range("K65").Value = 10000
$ws.Range("M65").Value = -6880

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2241.7778
$ws.Range("I81").Value = 2260.75
$ws.Range("K81").Value = 4521.5
$ws.Range("M81").Value = -3460.5

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 2241.7778
$ws.Range("I84").Value = 2260.75
$ws.Range("K84").Value = 22607.5
$ws.Range("M84").Value = -17303.5

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1649.8
$ws.Range("I122").Value = 1501.7142
$ws.Range("J122").Value = 1995.3334
$ws.Range("K122").Value = 4505.142599999999
$ws.Range("L122").Value = 5986.0002
$ws.Range("M122").Value = -2055.142599999999
$ws.Range("N122").Value = -10886.0002

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1951.7333
$ws.Range("I126").Value = 1606.4546
$ws.Range("K126").Value = 4819.3638
$ws.Range("M126").Value = -2349.3638

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2534.5
$ws.Range("I132").Value = 2534.5
$ws.Range("K132").Value = 7603.5
$ws.Range("M132").Value = -5073.5
